$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Fecha" column (B) shifts up by one appointment slot: each row's date
# becomes the date that used to belong to the row below it, and a new date
# ("Miércoles 05/06/2024") is appended for the freed-up last slot (B11),
# continuing the Mon/Tue/Wed appointment cadence.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r + 1, 2).Value2
}

$ws.Cells.Item(11, 2).Value = "Miércoles 05/06/2024"
